# Update iserv_stats for 2025-09 (row 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6287
$ws.Range("D22").Value = 5747670
$ws.Range("E22").Value = 914.2150469222205
$ws.Range("F22").Value = 8.228610776381483
$ws.Range("H22").Value = 24.99391087257197
